$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CustomerDetails")

# Update the country value in C5 from "Bangladesh" to "Afghanistan"
$ws.Range("C5").Value = '"Afghanistan"'

# Update the selected/active cell in the sheet view from D15 to F10
$ws.Range("F10").Select()
